$d = $word.ActiveDocument

# --- 1. Refresh the footer timestamp -------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
[void]$ftr.Range.Find.Execute("2025-06-30 12:12Z / ", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# --- 2. Add the regression-test character styles (b, i, sub, sup, u) -----
$b = $d.Styles.Add("b", 2)
$b.BaseStyle = "DefaultParagraphFont"
$b.Priority = 1
$b.QuickStyle = $true
$b.Font.Bold = $true

$i = $d.Styles.Add("i", 2)
$i.BaseStyle = "DefaultParagraphFont"
$i.Priority = 1
$i.QuickStyle = $true
$i.Font.Italic = $true

$sub = $d.Styles.Add("sub", 2)
$sub.BaseStyle = "DefaultParagraphFont"
$sub.Priority = 1
$sub.QuickStyle = $true
$sub.Font.Subscript = $true

$sup = $d.Styles.Add("sup", 2)
$sup.BaseStyle = "DefaultParagraphFont"
$sup.Priority = 1
$sup.QuickStyle = $true
$sup.Font.Superscript = $true

$u = $d.Styles.Add("u", 2)
$u.BaseStyle = "DefaultParagraphFont"
$u.Priority = 1
$u.QuickStyle = $true
$u.Font.Underline = 1
